# Add "getAllCacheNames" sheet to the JinZu ApiEngine test-data workbook
# (cache controller test case), mirroring the structure of the existing
# "getAllCacheStatsWithAuth" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet at the end of the workbook ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "getAllCacheNames"
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- 2. Fill in the cell values -------------------------------------------
# Write order matters: it controls the order new strings are appended to
# the shared-strings table, matching the authored file.
$json = "[`n    ""subscriptionCache"",`n    ""entityDataCache"",`n    ""graphQLSchemaCache"",`n    ""kgCache"",`n    ""connectorMetaDataCache""`n  ]"

$ws.Range("F2").Value = "Operate success."
$ws.Range("C2").Value = $json
$ws.Range("B2").Value = "get all cache names"
$ws.Range("A2").Value = "JinZu-ApiEngine-CacheController-Test-1"
$ws.Range("C1").Value = "rspData"

$ws.Range("A1").Value = "test-id"
$ws.Range("B1").Value = "description"
$ws.Range("D1").Value = "rspStatus"
$ws.Range("E1").Value = "rspCode"
$ws.Range("F1").Value = "rspMessage"

$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 0

# --- 3. Formatting: copy styles from the sibling "stats" sheet ------------
$stats = $wb.Worksheets.Item("getAllCacheStatsWithAuth")

$stats.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$stats.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)

$stats.Range("A2:B2").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)

$stats.Range("A2").Copy()
$ws.Range("D2:F2").PasteSpecial(-4122)

# wrap-text style used for the multi-line JSON payload cell
$leaseDetails = $wb.Worksheets.Item("getLeaseDetails")
$leaseDetails.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = $json

# --- 4. Column widths / row height ----------------------------------------
$ws.Columns.Item(1).ColumnWidth = 58.166666666666664
$ws.Columns.Item(2).ColumnWidth = 27.333333333333332
$ws.Columns.Item(3).ColumnWidth = 41.0
$ws.Columns.Item(6).ColumnWidth = 14.5

$ws.Rows.Item(2).RowHeight = 92.4

# --- 5. Selection / activation ---------------------------------------------
$ws.Range("A3").Select()
$newSheet.Activate()

Write-Output "getAllCacheNames sheet added"
